$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching style of existing headers (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the Save column data
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
